$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Convert the "Login" sub-story unit test into a "Register" sub-story unit test.
$ws.Range("E3").Value = "Brandon Goldwax"
$ws.Range("B4").Value = "[Sub-story] Register"
$ws.Range("B5").Value = "User is able to successfully register"
$ws.Range("A8").Value = "1- Navigate to the register page"
$ws.Range("A11").Value = "4 - Click on the register button"
$ws.Range("C8").Value = "The user is registered using the inputed username/email and password then is redirected to the login page"
$ws.Range("D8").Value = "The user cannot register if the username/email is already used"
$ws.Range("B13").Value = "Account details are registered in the database"
$ws.Range("B6").Value = "Hold a valid and unused username and password, "

$ws.Rows.Item(8).RowHeight = 57

$ws.Range("C8:C11").Select()
